# Auto update Excel log 2026-02-04 14:06:26
# Appends newly collected sensor readings to the PIR, Humidity and
# Temperature logs (Bathroom sensors), matching the master log format.
$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param($SheetName, $StartRow, $Rows)
    $ws = $wb.Worksheets.Item($SheetName)
    for ($i = 0; $i -lt $Rows.Count; $i++) {
        $r = $StartRow + $i
        $row = $Rows[$i]
        for ($c = 0; $c -lt $row.Count; $c++) {
            $val = $row[$c]
            $cell = $ws.Cells.Item($r, $c + 1)
            # Force plain text storage for tokens Excel would otherwise
            # auto-convert to a date serial or a percentage number,
            # so the logged value round-trips as the literal string.
            if ($val -match '^\d{4}-\d{2}-\d{2}$' -or $val -match '^-?\d+(\.\d+)?%$') {
                $cell.NumberFormat = "@"
            }
            $cell.Value = $val
        }
    }
}

# Date, Timestamp, Hour, Location, Value, Status
$pirRows = @(
    @("2026-02-04","14:05:22","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:05:24","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:05:30","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:05:35","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:05:40","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:05:45","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:05:50","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:05:55","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:06:00","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:06:05","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:06:10","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:06:15","14:00","Bathroom","No Motion","Inactive"),
    @("2026-02-04","14:06:20","14:00","Bathroom","No Motion","Inactive")
)

$humidityRows = @(
    @("2026-02-04","14:05:23","14:00","Bathroom","76.4%","Active"),
    @("2026-02-04","14:05:28","14:00","Bathroom","77.5%","Active"),
    @("2026-02-04","14:05:33","14:00","Bathroom","76.5%","Active"),
    @("2026-02-04","14:05:38","14:00","Bathroom","77.3%","Active"),
    @("2026-02-04","14:05:43","14:00","Bathroom","76.6%","Active"),
    @("2026-02-04","14:05:48","14:00","Bathroom","77.4%","Active"),
    @("2026-02-04","14:05:53","14:00","Bathroom","76.5%","Active"),
    @("2026-02-04","14:05:58","14:00","Bathroom","77.4%","Active"),
    @("2026-02-04","14:06:03","14:00","Bathroom","76.6%","Active"),
    @("2026-02-04","14:06:13","14:00","Bathroom","76.4%","Active"),
    @("2026-02-04","14:06:19","14:00","Bathroom","77.5%","Active")
)

$temperatureRows = @(
    @("2026-02-04","14:05:23","14:00","Bathroom","24.8C","Active"),
    @("2026-02-04","14:05:28","14:00","Bathroom","24.8C","Active"),
    @("2026-02-04","14:05:33","14:00","Bathroom","24.8C","Active"),
    @("2026-02-04","14:05:39","14:00","Bathroom","24.8C","Active"),
    @("2026-02-04","14:05:44","14:00","Bathroom","24.8C","Active"),
    @("2026-02-04","14:05:49","14:00","Bathroom","24.9C","Active"),
    @("2026-02-04","14:05:54","14:00","Bathroom","24.9C","Active"),
    @("2026-02-04","14:05:59","14:00","Bathroom","24.9C","Active"),
    @("2026-02-04","14:06:04","14:00","Bathroom","24.9C","Active"),
    @("2026-02-04","14:06:14","14:00","Bathroom","24.8C","Active"),
    @("2026-02-04","14:06:19","14:00","Bathroom","24.8C","Active")
)

Add-LogRows "PIR" 54 $pirRows
Add-LogRows "Humidity" 45 $humidityRows
Add-LogRows "Temperature" 45 $temperatureRows
